$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.563732385635376
$ws.Range("B1").Value = 2.643516302108765
$ws.Range("C1").Value = 6.484146118164062
$ws.Range("D1").Value = 1.910553097724915
$ws.Range("E1").Value = 1.601638793945312
